# Updated remaining queries for C3DC
# Fixes the JOIN conditions in the SQL queries (StatQuery + the 6 per-tab
# TabQuery cells) so they reference the renamed id columns
# (e.g. std.id -> std.study_id, prt.id -> prt.participant_id), matching the
# df_* dataframe column renames used elsewhere.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All the query cells that embed the old "<alias>.id" join conditions.
$queryCells = @("B2", "C2", "B3", "B4", "B5", "B6", "B7")

foreach ($cellAddr in $queryCells) {
    $cell = $ws.Range($cellAddr)
    $text = $cell.Value()
    if ($text -ne $null) {
        $text = $text -replace 'std\.id = prt\."study\.id"', 'std.study_id = prt."study.study_id"'
        $text = $text -replace 'prt\.id = dgn\."participant\.id"', 'prt.participant_id = dgn."participant.participant_id"'
        $text = $text -replace 'prt\.id = trt\."participant\.id"', 'prt.participant_id = trt."participant.participant_id"'
        $text = $text -replace 'prt\.id = trr\."participant\.id"', 'prt.participant_id = trr."participant.participant_id"'
        $text = $text -replace 'prt\.id = srv\."participant\.id"', 'prt.participant_id = srv."participant.participant_id"'
        $text = $text -replace 'std\.id = rfs\."study\.id"', 'std.study_id = rfs."study.study_id"'
        $cell.Value = $text
    }
}

# Widen column C to fit the updated query text (drop the stale "best fit"
# auto-size flag in favor of an explicit width) and move the active
# selection back up to B2.
$ws.Columns.Item(3).ColumnWidth = 67.67
$ws.Range("B2").Select() | Out-Null
